$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G column values from 55.2 to 55 for rows 3-8 and 24-26
$rows = @(3,4,5,6,7,8,24,25,26)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = 55
}

# Move the active cell selection to G27 (was H16)
$ws.Range("G27").Select()
